# Update the "cables" workbook ("Se ha actualizado el excel 'cables'"):
#
#  1. Split the thermal-breaker rating list out of the "cables" table into
#     its own sheet named "termic", inserted between "cables" and
#     "fusibles". The new sheet holds the standard list of breaker ratings
#     (A): 10, 16, 20, 25, 32, 40, 50, 63, 80, 100, 125, 160, 200, 250.
#
#  2. Remove the now-redundant first column ("Termic (A)") from the
#     "cables" sheet, so the wire-gauge table now starts with
#     "area (mm^2)" in column A (everything else shifts one column left).

$wb = $excel.ActiveWorkbook

$cables = $wb.Worksheets.Item("cables")
$fusibles = $wb.Worksheets.Item("fusibles")

# --- 1. Create the new "termic" sheet, positioned before "fusibles" ---
$termic = $wb.Worksheets.Add($fusibles)
$termic.Name = "termic"

$termic.Range("A1").Value = "Termic (A)"

$termicValues = @(10.0, 16.0, 20.0, 25.0, 32.0, 40.0, 50.0, 63.0, 80.0, 100.0, 125.0, 160.0, 200.0, 250.0)
$termicData = New-Object 'object[,]' $termicValues.Length,1
for ($i = 0; $i -lt $termicValues.Length; $i++) {
    $termicData[$i, 0] = $termicValues[$i]
}
$termic.Range("A2:A15").Value = $termicData

# --- 2. Drop the old "Termic (A)" column from "cables" (shifts table left) ---
$cables.Columns.Item(1).Delete()
